$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.201.15"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "3.114.01"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'579.80"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "'173.30"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("D9").Value = "'6.52"
$ws.Range("E9").Value = "  +1.28%  "
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").Value = "'36.90"
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("E14").Value = "  -1.63%  "
$ws.Range("D15").Value = "3.631.56"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").Value = "67.191.49"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("E17").Value = "  -1.42%  "
$ws.Range("D18").Value = "3.113.01"
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").Value = "'16.56"
$ws.Range("E19").Value = "  +1.64%  "
$ws.Range("D20").Value = "'490.92"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("E21").Value = "  +4.78%  "
$ws.Range("D22").Value = "'0.705"
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("D24").Value = "'13.21"
$ws.Range("E24").Value = "  -1.32%  "
$ws.Range("E25").Value = "  -3.12%  "
$ws.Range("D26").Value = "'10.58"
$ws.Range("E26").Value = "  +5.42%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "'7.96"
$ws.Range("E28").Value = "  -1.03%  "
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("D31").Value = "'28.42"
$ws.Range("E31").Value = "  -1.63%  "
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("D33").Value = "0.0₃0944"
$ws.Range("E33").Value = "  -6.34%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "'5.87"
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("D37").Value = "'47.12"
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("E38").Value = "  -3.68%  "
$ws.Range("E39").Value = "  -2.17%  "
$ws.Range("E40").Value = "  +0.89%  "
$ws.Range("D41").Value = "'8.49"
$ws.Range("E41").Value = "  -2.24%  "
$ws.Range("D42").Value = "'387.42"
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("D43").Value = "2.808.04"
$ws.Range("E43").Value = "  -1.52%  "
$ws.Range("E44").Value = "  -7.59%  "
$ws.Range("E45").Value = "  -2.46%  "
$ws.Range("D46").Value = "'135.17"
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").Value = "'24.98"
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("E49").Value = "  -0.93%  "
$ws.Range("E50").Value = "  -1.07%  "
$ws.Range("D51").Value = "'6.72"
$ws.Range("E51").Value = "  -2.31%  "
